# Update to downloadable table
$wb = $excel.ActiveWorkbook

$wsAnimal = $wb.Worksheets.Item("Animal Cases")
$wsHuman  = $wb.Worksheets.Item("Human Exposures")

# 1. Shrink the merged cell range on the Animal Cases sheet from B3:G7 to B3:E7
$wsAnimal.Range("B3:G7").UnMerge()
$wsAnimal.Range("B3:E7").Merge()

# 2. Column widths: Animal Cases A:F from 14 to 20
$wsAnimal.Range("A1:F1").ColumnWidth = 19.15

# 3. Column widths: Human Exposures A:H from 14 to 20
$wsHuman.Range("A1:H1").ColumnWidth = 19.15

# 4. Update the citation/footer text (shared by both sheets' B3 cell)
$newText = "This data is from the XXX. Please cite (YYY) if used. For more information, please contact Katie Hampson (Katie.Hampson@glasgow.ac.uk)"
$wsAnimal.Range("B3").Value = $newText
$wsHuman.Range("B3").Value = $newText

# 5. Font size for that caption text: 14 -> 12
$wsAnimal.Range("B3").Font.Size = 12
